$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: set a cell to a text value while preventing Excel's automatic
# type-inference (e.g. turning "2/26/2026" into a date serial, or "323"
# into a number). We temporarily force a Text number format, assign the
# value, then restore the cell style to Normal so no stray formatting is
# left behind on the cell.
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# Row 6
Set-TextValue $ws.Range("A6") "SNOW-712162"
Set-TextValue $ws.Range("B6") "2/26/2026"
Set-TextValue $ws.Range("C6") "local"
Set-TextValue $ws.Range("D6") "l@l.com"
Set-TextValue $ws.Range("E6") "323"
$ws.Range("F6").Value = 1
Set-TextValue $ws.Range("G6") "Family Ski Package"
$ws.Range("H6").Value = 32000
$ws.Range("I6").Value = 64000

# Row 7
Set-TextValue $ws.Range("A7") "SNOW-712162"
Set-TextValue $ws.Range("B7") "2/27/2026"
Set-TextValue $ws.Range("C7") "local"
Set-TextValue $ws.Range("D7") "l@l.com"
Set-TextValue $ws.Range("E7") "323"
$ws.Range("F7").Value = 1
Set-TextValue $ws.Range("G7") "Family Ski Package"
$ws.Range("H7").Value = 32000
$ws.Range("I7").Value = 64000
